$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8171977400779724
$ws.Range("B1").Value = 3.076165676116943
$ws.Range("C1").Value = 2.536687135696411
$ws.Range("D1").Value = 2.213939428329468
$ws.Range("E1").Value = 1.887593626976013
